# feat(logic), #55, #57: Add IFS function.
#
# Adds two demonstration rows (16 and 17) to the DATA sheet that exercise
# the new IFS() function, mirroring the existing IF()-based circular /
# non-circular example pairs already present in the sheet (e.g. rows 6/7).
# Row 16 is the "normal" (non-circular) case, row 17 is the circular case.
#
# Also enables iterative calculation (needed for the circular-reference
# examples on this "circular.xlsx" sheet) and moves the active selection,
# matching the authored workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn on iterative calculation (this sheet is full of intentional
# circular references) with a single iteration, and don't force a
# recalculation on save.
$excel.Iteration = $true
$excel.MaxIterations = 1
$excel.CalculateBeforeSave = $false

# Row 16: A16 = FALSE -> IFS() takes the non-circular TRUE branch.
$ws.Range("A16").Value = $false
$ws.Range("B16").Formula = "=IFS(NOT(A16),1,TRUE,C16)"
$ws.Range("C16").Formula = "=D16"
$ws.Range("D16").Formula = "=B16"
$ws.Range("E16").Formula = "=B16+1"

# Row 17: A17 = TRUE -> IFS() falls through to the circular branch.
$ws.Range("A17").Value = $true
$ws.Range("B17").Formula = "=IFS(NOT(A17),1,TRUE,C17)"
$ws.Range("C17").Formula = "=D17"
$ws.Range("D17").Formula = "=B17"
$ws.Range("E17").Formula = "=B17+1"

# Match the authored selection / active cell.
[void]$ws.Range("A15").Select()
